$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make "Alias Export" (column E) values unique by appending a
#     Revenue/Capex/Opex (+ share qualifier) suffix, per row. ---
$aliasMap = @{
    14 = "Total Amount Revenue"
    15 = "Relative Share in Percent Revenue NE"
    16 = "Absolute Share Revenue NE"
    17 = "Relative Share in Percent Revenue E"
    18 = "Absolute Share Revenue E"
    19 = "Relative Share in Percent Revenue NA"
    20 = "Absolute Share Revenue NA"
    21 = "Non-Aligned Activities Revenue"
    22 = "Relative Share in Percent Revenue"
    23 = "Absolute Share Revenue"
    24 = "Substantial Contribution to CCM In Percent - Eligible Revenue"
    25 = "Substantial Contribution to CCM In Percent - Aligned Revenue"
    26 = "Substantial Contribution to CCM In Percent - Of which use of proceeds Revenue"
    27 = "Substantial Contribution to CCM In Percent - Enabling Share Revenue"
    28 = "Substantial Contribution to CCM In Percent - Transitional Share Revenue"
    29 = "Substantial Contribution to CCA In Percent - Eligible Revenue"
    30 = "Substantial Contribution to CCA In Percent - Aligned Revenue"
    31 = "Substantial Contribution to CCA In Percent - Of which use of proceeds Revenue"
    32 = "Substantial Contribution to CCA In Percent - Enabling Share Revenue"
    33 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Eligible Revenue"
    34 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Aligned Revenue"
    35 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Of which use of proceeds Revenue"
    36 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Enabling Share Revenue"
    37 = "Substantial Contribution to CE In Percent - Eligible Revenue"
    38 = "Substantial Contribution to CE In Percent - Aligned Revenue"
    39 = "Substantial Contribution to CE In Percent - Of which use of proceeds Revenue"
    40 = "Substantial Contribution to CE In Percent - Enabling Share Revenue"
    41 = "Substantial Contribution to PPC In Percent - Eligible Revenue"
    42 = "Substantial Contribution to PPC In Percent - Aligned Revenue"
    43 = "Substantial Contribution to PPC In Percent - Of which use of proceeds Revenue"
    44 = "Substantial Contribution to PPC In Percent - Enabling Share Revenue"
    45 = "Substantial Contribution to BIO In Percent - Eligible Revenue"
    46 = "Substantial Contribution to BIO In Percent - Aligned Revenue"
    47 = "Substantial Contribution to BIO In Percent - Of which use of proceeds Revenue"
    48 = "Substantial Contribution to BIO In Percent - Enabling Share Revenue"
    49 = "Aligned Activities Revenue"
    50 = "Enabling Share In Percent Revenue"
    51 = "Transitional Share In Percent Revenue"
    52 = "Total Amount Capex"
    53 = "Relative Share in Percent Capex NE"
    54 = "Absolute Share Capex NE"
    55 = "Relative Share in Percent Capex E"
    56 = "Absolute Share Capex E"
    57 = "Relative Share in Percent Capex NA"
    58 = "Absolute Share Capex NA"
    59 = "Non-Aligned Activities Capex NA"
    60 = "Relative Share in Percent Capex A"
    61 = "Absolute Share Capex A"
    62 = "Substantial Contribution to CCM In Percent - Eligible Capex"
    63 = "Substantial Contribution to CCM In Percent - Aligned Capex"
    64 = "Substantial Contribution to CCM In Percent - Of which use of proceeds Capex"
    65 = "Substantial Contribution to CCM In Percent - Enabling Share Capex"
    66 = "Substantial Contribution to CCM In Percent - Transitional Share Capex"
    67 = "Substantial Contribution to CCA In Percent - Eligible Capex"
    68 = "Substantial Contribution to CCA In Percent - Aligned Capex"
    69 = "Substantial Contribution to CCA In Percent - Of which use of proceeds Capex"
    70 = "Substantial Contribution to CCA In Percent - Enabling Share Capex"
    71 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Eligible Capex"
    72 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Aligned Capex"
    73 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Of which use of proceeds Capex"
    74 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Enabling Share Capex"
    75 = "Substantial Contribution to CE In Percent - Eligible Capex"
    76 = "Substantial Contribution to CE In Percent - Aligned Capex"
    77 = "Substantial Contribution to CE In Percent - Of which use of proceeds Capex"
    78 = "Substantial Contribution to CE In Percent - Enabling Share Capex"
    79 = "Substantial Contribution to PPC In Percent - Eligible Capex"
    80 = "Substantial Contribution to PPC In Percent - Aligned Capex"
    81 = "Substantial Contribution to PPC In Percent - Of which use of proceeds Capex"
    82 = "Substantial Contribution to PPC In Percent - Enabling Share Capex"
    83 = "Substantial Contribution to BIO In Percent - Eligible Capex"
    84 = "Substantial Contribution to BIO In Percent - Aligned Capex"
    85 = "Substantial Contribution to BIO In Percent - Of which use of proceeds Capex"
    86 = "Substantial Contribution to BIO In Percent - Enabling Share Capex"
    87 = "Aligned Activities Capex"
    88 = "Enabling Share In Percent Capex"
    89 = "Transitional Share In Percent Capex"
    90 = "Total Amount Opex"
    91 = "Relative Share in Percent Opex NE"
    92 = "Absolute Share Opex NE"
    93 = "Relative Share in Percent Opex E"
    94 = "Absolute Share Opex E"
    95 = "Relative Share in Percent Opex NA"
    96 = "Absolute Share Opex NA"
    97 = "Non-Aligned Activities Opex"
    98 = "Relative Share in Percent Opex"
    99 = "Absolute Share Opex"
    100 = "Substantial Contribution to CCM In Percent - Aligned Opex"
    101 = "Substantial Contribution to CCA In Percent - Aligned Opex"
    102 = "Substantial Contribution to Sustainable Use and Protection of WTR In Percent - Aligned Opex"
    103 = "Substantial Contribution to CE In Percent - Aligned Opex"
    104 = "Substantial Contribution to PPC In Percent - Aligned Opex"
    105 = "Substantial Contribution to BIO In Percent - Aligned Opex"
    106 = "Aligned Activities Opex"
    107 = "Enabling Share In Percent Opex"
    108 = "Transitional Share In Percent Opex"
}

foreach ($rowNum in $aliasMap.Keys) {
    $ws.Cells.Item([int]$rowNum, 5).Value = $aliasMap[$rowNum]
}

# --- 2. Flag now-theoretically-unique column duplicates with a
#     conditional-formatting "duplicate values" rule (orange highlight),
#     matching the dxf (font FF9C5700 / fill FFFFEB9C) and the leftover
#     priority=2 numbering left behind by Excel UI workflow. ---
$dupRange = $ws.Range("E1:E1048576")

$scratchRule = $dupRange.FormatConditions.AddUniqueValues()
$scratchRule.DupeUnique = 0

$dupRule = $dupRange.FormatConditions.AddUniqueValues()
$dupRule.DupeUnique = 1
$dupRule.Font.Color = 22428
$dupRule.Interior.Color = 10284031

$scratchRule.Delete()

Write-Output "Applied unique alias export values and duplicate-highlight rule."
